$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("smzdm")

# Swap "active" login row from row 9 to row 14.
# Row 9: set H9 (login) back to "n", clear its mode/category/pages extras stay same text (unchanged) -- only H changes.
$ws.Range("H9").Value = "n"

# Row 14: set H14 (login) to "y", and copy mode/category/pages values that used to live on row 9.
$ws.Range("H14").Value = "y"
$ws.Range("I14").Value = "smzdm_share"
$ws.Range("J14").Value = "ele"
$ws.Range("K14").Value = "50-56"

# Update the selection to reflect the newly active row.
$ws.Range("K14").Select()
